# Applies the "worked on menu class" edit:
#  - Rename headers: "Order no." -> "item.no", "F_name" -> "name", "L_name" -> "price"
#  - Remove the old "Postcode"/"Phone" columns (D, E) from the header row
#  - Extend the numbered list in column A from 20 rows down to 28 rows (values 21-28 in rows 22-29)
#  - Move the active selection to B2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the now-unused header cells in columns D and E (Postcode, Phone)
$ws.Range("D1:E1").ClearContents() | Out-Null

# Update the remaining header labels
$ws.Range("A1").Value = "item.no"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "price"

# Extend the numbered list (column A) with rows 22-29 holding values 21-28
for ($row = 22; $row -le 29; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1
}

# Update the selected cell shown when the workbook is opened
$ws.Range("B2").Select() | Out-Null
